$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 12:20"

$ws.Range("B4").Value = 142746
$ws.Range("C4").Value = 286
$ws.Range("E4").Value = 135695

$ws.Range("A9").Value = "Iran"
$ws.Range("B9").Value = 41495
$ws.Range("C9").Value = 3186
$ws.Range("D9").Value = 13911
$ws.Range("E9").Value = 24827
$ws.Range("F9").Value = 3511
$ws.Range("G9").Value = 117
$ws.Range("H9").Value = 2757

$ws.Range("A10").Value = "Francia"
$ws.Range("B10").Value = 40174
$ws.Range("D10").Value = 7202
$ws.Range("E10").Value = 30366
$ws.Range("F10").Value = 4632
$ws.Range("H10").Value = 2606

$ws.Range("A13").Value = "Belgica"
$ws.Range("B13").Value = 11899
$ws.Range("C13").Value = 1063
$ws.Range("D13").Value = 1527
$ws.Range("E13").Value = 9859
$ws.Range("F13").Value = 927
$ws.Range("G13").Value = 82
$ws.Range("H13").Value = 513

$ws.Range("A14").Value = "Paises Bajos"
$ws.Range("B14").Value = 10866
$ws.Range("D14").Value = 250
$ws.Range("E14").Value = 9845
$ws.Range("F14").Value = 972
$ws.Range("H14").Value = 771

$ws.Range("B17").Value = 9103
$ws.Range("C17").Value = 315
$ws.Range("E17").Value = 8538

$ws.Range("E20").Value = 4197
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 16

$ws.Range("B21").Value = 4313
$ws.Range("C21").Value = 29
$ws.Range("E21").Value = 4277
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 29

$ws.Range("A30").Value = "Rumania"
$ws.Range("B30").Value = 1952
$ws.Range("C30").Value = 137
$ws.Range("D30").Value = 206
$ws.Range("E30").Value = 1702
$ws.Range("F30").Value = 31
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 44

$ws.Range("A31").Value = "Luxemburgo"
$ws.Range("B31").Value = 1950
$ws.Range("D31").Value = 40
$ws.Range("E31").Value = 1889
$ws.Range("F31").Value = 25
$ws.Range("H31").Value = 21

$ws.Range("A32").Value = "Ecuador"
$ws.Range("B32").Value = 1924
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 3
$ws.Range("E32").Value = 1863
$ws.Range("F32").Value = 58
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 58

$ws.Range("A33").Value = "Polonia"
$ws.Range("B33").Value = 1905
$ws.Range("C33").Value = 43
$ws.Range("D33").Value = 7
$ws.Range("E33").Value = 1872
$ws.Range("F33").Value = 3
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = 26

$ws.Range("A34").Value = "Japon"
$ws.Range("B34").Value = 1866
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 424
$ws.Range("E34").Value = 1388
$ws.Range("F34").Value = 56
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 54

$ws.Range("A35").Value = "Rusia"
$ws.Range("B35").Value = 1836
$ws.Range("C35").Value = 302
$ws.Range("D35").Value = 66
$ws.Range("E35").Value = 1761
$ws.Range("F35").Value = 8
$ws.Range("H35").Value = 9

$ws.Range("E38").Value = 1286
$ws.Range("G38").Value = 2
$ws.Range("H38").Value = 9

$ws.Range("A40").Value = "Finlandia"
$ws.Range("B40").Value = 1343
$ws.Range("C40").Value = 103
$ws.Range("D40").Value = 10
$ws.Range("E40").Value = 1322
$ws.Range("F40").Value = 32
$ws.Range("H40").Value = 11

$ws.Range("A41").Value = "Arabia Saudita"
$ws.Range("B41").Value = 1299
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 66
$ws.Range("E41").Value = 1225
$ws.Range("F41").Value = 12
$ws.Range("H41").Value = 8

$ws.Range("E53").Value = 685
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 14

$ws.Range("A95").Value = "Senegal"
$ws.Range("B95").Value = 162
$ws.Range("C95").Value = 20
$ws.Range("D95").Value = 27
$ws.Range("E95").Value = 135
$ws.Range("F95").Value = 0

$ws.Range("A96").Value = "Islas Feroe"
$ws.Range("B96").Value = 159
$ws.Range("D96").Value = 70
$ws.Range("E96").Value = 89
$ws.Range("H96").Value = 0

$ws.Range("A97").Value = "Ghana"
$ws.Range("B97").Value = 152
$ws.Range("E97").Value = 145
$ws.Range("F97").Value = 1
$ws.Range("H97").Value = 5

$ws.Range("A98").Value = "Malta"
$ws.Range("B98").Value = 151
$ws.Range("D98").Value = 2
$ws.Range("E98").Value = 149
$ws.Range("F98").Value = 4
$ws.Range("H98").Value = 0

$ws.Range("A99").Value = "Uzbekistan"
$ws.Range("B99").Value = 144
$ws.Range("D99").Value = 7
$ws.Range("E99").Value = 135
$ws.Range("F99").Value = 8
$ws.Range("H99").Value = 2

$ws.Range("E148").Value = 15
$ws.Range("G148").Value = 2
$ws.Range("H148").Value = 3
